$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT, even when it looks like a number
# (e.g. a bare distance reading such as "185.04..."). The logger always wrote
# these columns as text, never as numeric cells, so a numeric-looking value
# must not be auto-converted to a real number by Excel.
# Briefly mark the cell as Text (@) so the value is accepted verbatim, then
# restore the cell to the workbooks normal (General) style.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Bug: the logger crashed if a trial started with the cursor already at the
# target center, truncating the captured-trials table. With the fix, a 4th
# trial now makes it into the log, and the table is rewritten with the
# corrected run data.

# Make room for the new (4th) trial row by pushing the "Subject Code:" row
# down from row 5 to row 6.
$ws.Rows.Item(5).Insert()

# Write the corrected trial data into rows 2-5.

# Row 2
$ws.Range('A2').Value = '(564, 367)'
$ws.Range('B2').Value = '(379, 371)'
Set-TextValue 'C2' '185.04323819042943'
$ws.Range('D2').Value = '0:00:01.225507'
$ws.Range('E2').Value = '[[479, 367, datetime.timedelta(microseconds=116545), 416.9553640327558], [427, 362, datetime.timedelta(microseconds=233742), 105.12474828050938, -1334.0803781615905], [387, 367, datetime.timedelta(microseconds=347788), 54.51952285715116, -145.50595599433626], [383, 360, datetime.timedelta(microseconds=460534), 8.234456442061992, -100.50303867920537], [379, 371, datetime.timedelta(microseconds=568476), 9.684743123859638, 2.551183659112515], [381, 367, datetime.timedelta(microseconds=684748), 3.072021014299863, -9.657161626700297], [381, 373, datetime.timedelta(microseconds=794539), 3.5520247869798998, 0.6041286490405591], [379, 369, datetime.timedelta(microseconds=910617), 2.3100384085733103, -1.3638954449637877], [381, 365, datetime.timedelta(seconds=1, microseconds=17284), 121.70563790209458, 6907.8685196436745], [381, 371, datetime.timedelta(seconds=1, microseconds=116432), 24.23923167361398, -837.1101263267882]]'
$ws.Range('F2').Value = '2022-07-18 14:08:43.849793'
Set-TextValue 'G2' '87.03885648216496'
$ws.Range('H2').Value = 'hit'
$ws.Range('I2').Value = '(416.0, 360.0)'
Set-TextValue 'J2' '38.600518131237564'

# Row 3
$ws.Range('A3').Value = '(576, 358)'
$ws.Range('B3').Value = '(416, 369)'
Set-TextValue 'C3' '160.37767924496225'
$ws.Range('D3').Value = '0:00:00.466837'
$ws.Range('E3').Value = '[[516, 354, datetime.timedelta(microseconds=119004), 273.7557543048282], [458, 367, datetime.timedelta(microseconds=231039), 121.01145404937245, -661.1191195229193], [418, 367, datetime.timedelta(microseconds=344831), 54.5624227949773, -192.70028290494517]]'
$ws.Range('F3').Value = '2022-07-18 14:08:51.289994'
Set-TextValue 'G3' '75.43690838559337'
$ws.Range('H3').Value = 'hit'
$ws.Range('I3').Value = '(416.0, 360.0)'
Set-TextValue 'J3' '9.0'

# Row 4
$ws.Range('A4').Value = '(593, 377)'
$ws.Range('B4').Value = '(454, 371)'
Set-TextValue 'C4' '139.12943613772032'
$ws.Range('D4').Value = '0:00:00.562122'
$ws.Range('E4').Value = '[[555, 375, datetime.timedelta(microseconds=122082), 104.02843990105013], [485, 371, datetime.timedelta(microseconds=232791), 141.6705917678617, 161.69934347466858], [458, 373, datetime.timedelta(microseconds=343206), 37.105396134542076, -304.671817023361], [454, 373, datetime.timedelta(microseconds=452708), 4.156059715051383, -72.78275714034365]]'
$ws.Range('F4').Value = '2022-07-18 14:08:57.435283'
Set-TextValue 'G4' '65.44236440552031'
$ws.Range('H4').Value = 'hit'
$ws.Range('I4').Value = '(416.0, 360.0)'
Set-TextValue 'J4' '39.56008088970496'

# Row 5
$ws.Range('A5').Value = '(582, 375)'
$ws.Range('B5').Value = '(437, 369)'
Set-TextValue 'C5' '145.1240848377691'
$ws.Range('D5').Value = '0:00:00.539133'
$ws.Range('E5').Value = '[[512, 369, datetime.timedelta(microseconds=95883), 344.6560595914215], [468, 367, datetime.timedelta(microseconds=205065), 101.02975025206085, -1188.0443241867733], [441, 373, datetime.timedelta(microseconds=321080), 40.518878855951336, -188.46041919804884], [437, 367, datetime.timedelta(microseconds=424116), 7.997550145853295, -76.68026839378388]]'
$ws.Range('F5').Value = '2022-07-18 14:09:03.708945'
Set-TextValue 'G5' '68.2620695348025'
$ws.Range('H5').Value = 'hit'
$ws.Range('I5').Value = '(416.0, 360.0)'
Set-TextValue 'J5' '22.847319317591726'
